$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update Row 4, column D (Exceute) from "Y" to "N"
$ws.Range("D4").Value = "N"

# Add new Row 5 with the new test case data
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Create Datasets from connectors"
$ws.Range("C5").Value = "Create Datasets from connectors"
$ws.Range("D5").Value = "Y"

# Update the active selection to match the target state
$ws.Range("D11").Select()
